function Compare-Ordinal($strA, $strB) {
    $lenA = $strA.Length
    $lenB = $strB.Length
    $minLen = [Math]::Min($lenA, $lenB)
    for ($pos = 0; $pos -lt $minLen; $pos++) {
        $codeA = [int]$strA[$pos]
        $codeB = [int]$strB[$pos]
        if ($codeA -lt $codeB) { return -1 }
        if ($codeA -gt $codeB) { return 1 }
    }
    if ($lenA -lt $lenB) { return -1 }
    if ($lenA -gt $lenB) { return 1 }
    return 0
}

function Sort-Ordinal($arrIn) {
    $list = @($arrIn)
    $cnt = $list.Count
    for ($idx = 1; $idx -lt $cnt; $idx++) {
        $key = $list[$idx]
        $jdx = $idx - 1
        while ($jdx -ge 0 -and (Compare-Ordinal $list[$jdx] $key) -gt 0) {
            $list[$jdx+1] = $list[$jdx]
            $jdx = $jdx - 1
        }
        $list[$jdx+1] = $key
    }
    return $list
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRowUsed = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row

for ($rowNum = 2; $rowNum -le $lastRowUsed; $rowNum++) {
    $gCell = $ws.Cells.Item($rowNum, 7)
    $origVal = $gCell.Value()
    if ($origVal -ne $null -and $origVal -ne "") {
        $nameParts = $origVal -split ", "
        $sortedParts = Sort-Ordinal $nameParts
        $joinedVal = [string]::Join(", ", $sortedParts)
        if ($joinedVal -ne $origVal) {
            $gCell.Value = $joinedVal
        }
    }
}
